$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.867.74'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '2.296.06'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.95'
$ws.Range("E5").Value = '  +17.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '270.18'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.17'
$ws.Range("E10").Value = '  +6.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0945'
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.05'
$ws.Range("E12").Value = '  +14.81%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.91'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '2.639.62'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '2.295.48'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").Value = '43.745.15'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.86'
$ws.Range("E20").Value = '  +10.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.22'
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.05'
$ws.Range("E23").Value = '  +12.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.02'
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.67'
$ws.Range("E25").Value = '  +5.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.63'
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.88'
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '41.77'
$ws.Range("E29").Value = '  +8.16%  '
$ws.Range("B30").Value = 'WEMIXToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.38'
$ws.Range("E30").Value = '  -2.01%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.27'
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '175.72'
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.58'
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0927'
$ws.Range("E34").Value = '  +3.63%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").Value = '  +5.22%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.127'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.67'
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0365'
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.107'
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.84'
$ws.Range("E40").Value = '  +11.48%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.91'
$ws.Range("E41").Value = '  +13.29%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.13'
$ws.Range("E42").Value = '  +14.99%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.242'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("E44").Value = '  +3.16%  '
$ws.Range("B45").Value = 'THORChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.29'
$ws.Range("E45").Value = '  +21.03%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.39'
$ws.Range("E47").Value = '  +3.10%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.82'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0997'
$ws.Range("E49").Value = '  -2.65%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.37'
$ws.Range("E50").Value = '  +3.09%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.24'
$ws.Range("E51").Value = '  +3.00%  '
